# Applies the v1.1.0 data update described in the commit "Added 1.1.0 version of data".
#
# The underlying dataset refresh reassigned the per-council metric rows (columns
# E:S - scores, weighted_total, quintile, pop_bucket, region, ruc_cluster,
# political_control) to different councils (column C/D) within several row
# blocks (rows 31-32, 44-51, 123-124, 146-147 on Sheet1). Each row below is
# rewritten in full (official-name, local-authority-code and every score/
# attribute column) to its new, authoritative value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C31").Value = 'Babergh District Council'
$ws.Range("D31").Value = 'BAB'
$ws.Range("E31").Value = 0.6666666666666666
$ws.Range("F31").Value = 0.7222222222222222
$ws.Range("G31").Value = 0.7142857142857143
$ws.Range("H31").Value = 0.6666666666666666
$ws.Range("I31").Value = 0.4
$ws.Range("J31").Value = 0.5
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0.6
$ws.Range("M31").Value = 0.75
$ws.Range("N31").Value = 0.5779761904761905
$ws.Range("O31").Value = 3
$ws.Range("P31").Value = '80k - 90k'
$ws.Range("Q31").Value = 'East of England'
$ws.Range("R31").Value = 'Rural'
$ws.Range("S31").Value = 'Coalition/Minority'
$ws.Range("C32").Value = 'Mid Suffolk District Council'
$ws.Range("D32").Value = 'MSU'
$ws.Range("E32").Value = 0.6666666666666666
$ws.Range("F32").Value = 0.7222222222222222
$ws.Range("G32").Value = 0.7142857142857143
$ws.Range("H32").Value = 0.6666666666666666
$ws.Range("I32").Value = 0.4
$ws.Range("J32").Value = 0.5
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0.6
$ws.Range("M32").Value = 0.75
$ws.Range("N32").Value = 0.5779761904761905
$ws.Range("O32").Value = 4
$ws.Range("P32").Value = '90k - 100k'
$ws.Range("Q32").Value = 'East of England'
$ws.Range("R32").Value = 'Sparse and rural'
$ws.Range("S32").Value = 'Coalition/Minority'
$ws.Range("C44").Value = 'Adur District Council'
$ws.Range("D44").Value = 'ADU'
$ws.Range("E44").Value = 0.6190476190476191
$ws.Range("F44").Value = 0.6666666666666666
$ws.Range("G44").Value = 0.2857142857142857
$ws.Range("H44").Value = 0.5555555555555556
$ws.Range("I44").Value = 0.6
$ws.Range("J44").Value = 0.75
$ws.Range("K44").Value = 0.4
$ws.Range("L44").Value = 0.6
$ws.Range("M44").Value = 0.5
$ws.Range("N44").Value = 0.5415476190476191
$ws.Range("O44").Value = 2
$ws.Range("P44").Value = '50k - 80k'
$ws.Range("Q44").Value = 'South East'
$ws.Range("R44").Value = 'Urban'
$ws.Range("S44").Value = 'Conservative'
$ws.Range("C45").Value = 'Basingstoke and Deane Borough Council'
$ws.Range("D45").Value = 'BAN'
$ws.Range("E45").Value = 0.6190476190476191
$ws.Range("F45").Value = 0.4444444444444444
$ws.Range("G45").Value = 0.8571428571428571
$ws.Range("H45").Value = 0.5555555555555556
$ws.Range("I45").Value = 0.8
$ws.Range("J45").Value = 0.5
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0.4
$ws.Range("M45").Value = 0.5
$ws.Range("N45").Value = 0.5414285714285715
$ws.Range("O45").Value = 4
$ws.Range("P45").Value = '170k - 250k'
$ws.Range("Q45").Value = 'South East'
$ws.Range("R45").Value = 'Urban with rural areas'
$ws.Range("S45").Value = 'Conservative'
$ws.Range("C46").Value = 'Carlisle City Council'
$ws.Range("D46").Value = 'CAR'
$ws.Range("E46").Value = 0.6666666666666666
$ws.Range("F46").Value = 0.2777777777777778
$ws.Range("G46").Value = 0.7142857142857143
$ws.Range("H46").Value = 0.6666666666666666
$ws.Range("I46").Value = 1
$ws.Range("J46").Value = 0.75
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0.4
$ws.Range("M46").Value = 0.25
$ws.Range("N46").Value = 0.5388095238095237
$ws.Range("O46").Value = 2
$ws.Range("P46").Value = '110k - 120k'
$ws.Range("Q46").Value = 'North West'
$ws.Range("R46").Value = 'Urban with rural areas'
$ws.Range("S46").Value = 'Conservative'
$ws.Range("C47").Value = 'West Lancashire Borough Council'
$ws.Range("D47").Value = 'WLA'
$ws.Range("E47").Value = 0.6190476190476191
$ws.Range("F47").Value = 0.5
$ws.Range("G47").Value = 0.8571428571428571
$ws.Range("H47").Value = 0.5555555555555556
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = 0.5
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0.2
$ws.Range("M47").Value = 0.25
$ws.Range("N47").Value = 0.5372619047619047
$ws.Range("O47").Value = 2
$ws.Range("P47").Value = '110k - 120k'
$ws.Range("Q47").Value = 'North West'
$ws.Range("R47").Value = 'Urban with rural areas'
$ws.Range("S47").Value = 'Coalition/Minority'
$ws.Range("C48").Value = 'St Albans City and District Council'
$ws.Range("D48").Value = 'SAL'
$ws.Range("E48").Value = 0.5714285714285714
$ws.Range("F48").Value = 0.6666666666666666
$ws.Range("G48").Value = 0.4285714285714285
$ws.Range("H48").Value = 0.7777777777777778
$ws.Range("I48").Value = 0.6
$ws.Range("J48").Value = 0.5
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0.4
$ws.Range("M48").Value = 0.75
$ws.Range("N48").Value = 0.5291666666666667
$ws.Range("O48").Value = 5
$ws.Range("P48").Value = '140k - 160k'
$ws.Range("Q48").Value = 'East of England'
$ws.Range("R48").Value = 'Urban'
$ws.Range("S48").Value = 'Liberal Democrat'
$ws.Range("C49").Value = 'Worthing Borough Council'
$ws.Range("D49").Value = 'WOT'
$ws.Range("E49").Value = 0.6190476190476191
$ws.Range("F49").Value = 0.6666666666666666
$ws.Range("G49").Value = 0.2857142857142857
$ws.Range("H49").Value = 0.5555555555555556
$ws.Range("I49").Value = 0.6
$ws.Range("J49").Value = 0.75
$ws.Range("K49").Value = 0.4
$ws.Range("L49").Value = 0.6
$ws.Range("M49").Value = 0.25
$ws.Range("N49").Value = 0.529047619047619
$ws.Range("O49").Value = 3
$ws.Range("P49").Value = '110k - 120k'
$ws.Range("Q49").Value = 'South East'
$ws.Range("R49").Value = 'Urban'
$ws.Range("S49").Value = 'Conservative'
$ws.Range("C50").Value = 'Warwick District Council'
$ws.Range("D50").Value = 'WAW'
$ws.Range("E50").Value = 0.9047619047619048
$ws.Range("F50").Value = 0.5555555555555556
$ws.Range("G50").Value = 0.4285714285714285
$ws.Range("H50").Value = 0.6666666666666666
$ws.Range("I50").Value = 0.8
$ws.Range("J50").Value = 0.5
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0.75
$ws.Range("N50").Value = 0.5258333333333334
$ws.Range("O50").Value = 4
$ws.Range("P50").Value = '140k - 160k'
$ws.Range("Q50").Value = 'West Midlands'
$ws.Range("R50").Value = 'Urban'
$ws.Range("S50").Value = 'Coalition/Minority'
$ws.Range("C51").Value = 'Colchester Borough Council'
$ws.Range("D51").Value = 'COL'
$ws.Range("E51").Value = 0.4761904761904762
$ws.Range("F51").Value = 0.2777777777777778
$ws.Range("G51").Value = 0.8571428571428571
$ws.Range("H51").Value = 0.6666666666666666
$ws.Range("I51").Value = 0.6
$ws.Range("J51").Value = 1
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0.4
$ws.Range("M51").Value = 0.5
$ws.Range("N51").Value = 0.5166666666666666
$ws.Range("O51").Value = 3
$ws.Range("P51").Value = '170k - 250k'
$ws.Range("Q51").Value = 'East of England'
$ws.Range("R51").Value = 'Urban with rural areas'
$ws.Range("S51").Value = 'Coalition/Minority'
$ws.Range("C123").Value = 'Test Valley Borough Council'
$ws.Range("D123").Value = 'TES'
$ws.Range("E123").Value = 0.3333333333333333
$ws.Range("F123").Value = 0.05555555555555555
$ws.Range("G123").Value = 0.2857142857142857
$ws.Range("H123").Value = 0.4444444444444444
$ws.Range("I123").Value = 0.4
$ws.Range("J123").Value = 0.25
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = 0
$ws.Range("N123").Value = 0.2203571428571429
$ws.Range("O123").Value = 4
$ws.Range("P123").Value = '120k - 130k'
$ws.Range("Q123").Value = 'South East'
$ws.Range("R123").Value = 'Urban with rural areas'
$ws.Range("S123").Value = 'Conservative'
$ws.Range("C124").Value = 'Hambleton District Council'
$ws.Range("D124").Value = 'HAE'
$ws.Range("E124").Value = 0.1904761904761905
$ws.Range("F124").Value = 0.1666666666666667
$ws.Range("G124").Value = 0.4285714285714285
$ws.Range("H124").Value = 0.3333333333333333
$ws.Range("I124").Value = 0.2
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0.2
$ws.Range("M124").Value = 0.25
$ws.Range("N124").Value = 0.2203571428571429
$ws.Range("O124").Value = 4
$ws.Range("P124").Value = '80k - 90k'
$ws.Range("Q124").Value = 'Yorkshire and The Humber'
$ws.Range("R124").Value = 'Rural'
$ws.Range("S124").Value = 'Conservative'
$ws.Range("C146").Value = 'Vale of White Horse District Council'
$ws.Range("D146").Value = 'VAL'
$ws.Range("E146").Value = 0
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 0
$ws.Range("I146").Value = 0
$ws.Range("J146").Value = 0
$ws.Range("K146").Value = 0
$ws.Range("L146").Value = 0
$ws.Range("M146").Value = 0
$ws.Range("N146").Value = 0
$ws.Range("O146").Value = 5
$ws.Range("P146").Value = '140k - 160k'
$ws.Range("Q146").Value = 'South East'
$ws.Range("R146").Value = 'Urban with rural areas'
$ws.Range("S146").Value = 'Liberal Democrat'
$ws.Range("C147").Value = 'Tamworth Borough Council'
$ws.Range("D147").Value = 'TAW'
$ws.Range("E147").Value = 0
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 0
$ws.Range("I147").Value = 0
$ws.Range("J147").Value = 0
$ws.Range("K147").Value = 0
$ws.Range("L147").Value = 0
$ws.Range("M147").Value = 0
$ws.Range("N147").Value = 0
$ws.Range("O147").Value = 2
$ws.Range("P147").Value = '50k - 80k'
$ws.Range("Q147").Value = 'West Midlands'
$ws.Range("R147").Value = 'Urban'
$ws.Range("S147").Value = 'Conservative'
